# Acceptance Test Plan - record Sprint 1 results for the "hint" user story
# acceptance criteria, and try to add a hint button (per commit message the
# author could not find the actual source file to wire up the button, so
# only the test-plan bookkeeping changed).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Plan")

# Rows that passed Sprint-1 testing (Status column C) -- this also covers
# the already-existing rows 2-34 plus row 36.
$passRows = 2..34 + 36
foreach ($r in $passRows) {
    $ws.Cells.Item($r, 3).Value = "Pass"
}

# Rows that failed Sprint-1 testing (Status column C) for the new "hint"
# stories.
$failRows = @(35, 37, 38, 39, 40, 41, 42)
foreach ($r in $failRows) {
    $ws.Cells.Item($r, 3).Value = "Fail"
}

# Comments column (D) for the failed rows -- "Not implemented" is used for
# most of them.
$notImplementedRows = @(37, 38, 39, 41, 42)
foreach ($r in $notImplementedRows) {
    $ws.Cells.Item($r, 4).Value = "Not implemented - MA, RT, ST 11/10/2021"
}

# Row 40 has its own, more specific, failure comment.
$ws.Cells.Item(40, 4).Value = "Just logs back into the original game. - MA, RT, ST 11/10/2021"

# Row 35 has its own failure comment too -- added last.
$ws.Cells.Item(35, 4).Value = "The page just refreshes when pressing the hoem button in the game. - MA, RT, ST 11/10/2021"

# Leave the sheet scrolled/selected where the author ended up working.
$ws.Activate()
$ws.Range("C9").Select()
